# Arestides_resume.docx edit:
#   1. "WORK EX" + bookmark + "PERIENCE" -> single run "WORK EXPERIENCE"
#      (the stray _GoBack bookmark that used to sit inside the heading
#       goes away as a natural side effect of the text being replaced).
#   2. "Web development: Capable with HTML, CSS, and MySQL" bullet gets
#      "/ mobile" and "Android development," worked into it, typed as two
#      separate edits (so the runs stay split exactly like real Word
#      leaves them), and the _GoBack bookmark ends up marking the spot
#      right after the last thing that was typed.

$d = $word.ActiveDocument

# --- 1. Heading: "WORK EX" / bookmark / "PERIENCE" -> "WORK EXPERIENCE" ---
# Find/Replace ignores the (empty) bookmark sitting between the two runs,
# matches across it, and rewrites the whole heading as one clean run -
# which also removes the now-enclosed _GoBack bookmark.
$headingFound = $d.Content.Find.Execute(
    "WORK EXPERIENCE", $false, $false, $false, $false, $false,
    $true, 1, $false, "WORK EXPERIENCE", 2)

# --- 2. Tech-skills bullet: weave in "/ mobile" and "Android development," ---
$bullet = $d.Content
$bulletFound = $bullet.Find.Execute(
    "Web development: Capable with HTML, CSS, and MySQL", $false, $false,
    $false, $false, $false, $true, 1, $false, "", 0)
$bulletStart = $bullet.Start

# Track the edits as separate insertions so they don't get silently
# re-coalesced back into the original single run.
$d.TrackRevisions = $true

# Type " / mobile" right after "Web".
$afterWeb = $d.Range($bulletStart + 3, $bulletStart + 3)
$afterWeb.InsertAfter(" / mobile")

# Type " Android development," right after "...HTML, CSS,".
$prefixLen = ("Web / mobile development: Capable with HTML, CSS,").Length
$afterCss = $d.Range($bulletStart + $prefixLen, $bulletStart + $prefixLen)
$afterCss.InsertAfter(" Android development,")

$d.TrackRevisions = $false
$d.AcceptAllRevisions()

# --- 3. Re-drop the _GoBack bookmark where the typing left off ---
# (right before the trailing " and MySQL")
$finalPrefixLen = ("Web / mobile development: Capable with HTML, CSS, Android development,").Length
$goBackPos = $bulletStart + $finalPrefixLen
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

Write-Output "headingFound=$headingFound bulletFound=$bulletFound"
